$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "43.888.77"
Set-TextValue "E2" "  -0.65%  "
Set-TextValue "D3" "2.349.84"
Set-TextValue "E3" "  -0.44%  "
Set-TextValue "E4" "  +0.23%  "
Set-TextValue "D5" "0.676"
Set-TextValue "E5" "  -0.16%  "
Set-TextValue "D6" "238.67"
Set-TextValue "E6" "  +0.05%  "
Set-TextValue "D7" "73.43"
Set-TextValue "E7" "  +0.61%  "
Set-TextValue "E8" "  +0.00%  "
Set-TextValue "D9" "0.591"
Set-TextValue "E9" "  +8.51%  "
Set-TextValue "D10" "0.0999"
Set-TextValue "E10" "  -0.66%  "
Set-TextValue "D11" "57.26"
Set-TextValue "E11" "  +0.03%  "
Set-TextValue "D12" "32.10"
Set-TextValue "E12" "  +9.00%  "
Set-TextValue "D13" "7.26"
Set-TextValue "E13" "  +8.23%  "
Set-TextValue "E14" "  +0.17%  "
Set-TextValue "D15" "2.703.75"
Set-TextValue "E15" "  -0.34%  "
Set-TextValue "D16" "16.51"
Set-TextValue "E16" "  -1.90%  "
Set-TextValue "D17" "0.895"
Set-TextValue "E17" "  -1.24%  "
Set-TextValue "D18" "2.363.76"
Set-TextValue "E18" "  +0.06%  "
Set-TextValue "D19" "43.816.34"
Set-TextValue "E19" "  -0.56%  "
Set-TextValue "E20" "  -1.41%  "
Set-TextValue "D21" "6.69"
Set-TextValue "E21" "  +3.49%  "
Set-TextValue "D22" "76.62"
Set-TextValue "E22" "  -1.73%  "
Set-TextValue "D23" "257.06"
Set-TextValue "E23" "  +0.52%  "
Set-TextValue "E24" "  +21.45%  "
Set-TextValue "E26" "  -2.90%  "
Set-TextValue "E27" "  -1.53%  "
Set-TextValue "D28" "10.69"
Set-TextValue "E28" "  +1.76%  "
Set-TextValue "D29" "2.27"
Set-TextValue "E29" "  +1.73%  "
Set-TextValue "D30" "22.57"
Set-TextValue "E30" "  +0.53%  "
Set-TextValue "D31" "175.54"
Set-TextValue "E31" "  +1.51%  "
Set-TextValue "D32" "0.128"
Set-TextValue "E32" "  -3.10%  "
Set-TextValue "E33" "  +2.59%  "
Set-TextValue "D34" "0.0762"
Set-TextValue "E34" "  +3.95%  "
Set-TextValue "D35" "5.19"
Set-TextValue "E35" "  -0.26%  "
Set-TextValue "D36" "5.42"
Set-TextValue "E36" "  +3.29%  "
Set-TextValue "D37" "3.74"
Set-TextValue "E37" "  -4.93%  "
Set-TextValue "D38" "2.34"
Set-TextValue "E38" "  -4.08%  "
Set-TextValue "D39" "6.26"
Set-TextValue "E39" "  -2.97%  "
Set-TextValue "D40" "0.0277"
Set-TextValue "E40" "  +2.23%  "
Set-TextValue "D41" "0.109"
Set-TextValue "E41" "  +11.29%  "
Set-TextValue "D42" "0.205"
Set-TextValue "E42" "  +12.62%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D43" "9.00"
Set-TextValue "E43" "  +1.61%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D44" "18.89"
Set-TextValue "E44" "  -3.69%  "
Set-TextValue "E45" "  +0.10%  "
Set-TextValue "D46" "4.75"
Set-TextValue "E46" "  +6.16%  "
Set-TextValue "D47" "57.83"
Set-TextValue "E47" "  +9.49%  "
Set-TextValue "D48" "2.50"
Set-TextValue "E48" "  +6.16%  "
Set-TextValue "D49" "1.23"
Set-TextValue "E49" "  -1.35%  "
Set-TextValue "E50" "  -0.56%  "
Set-TextValue "D51" "99.59"
Set-TextValue "E51" "  +0.96%  "
